$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '28.139.37'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +2.82%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.777.35'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -0.68%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.003'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.07%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '338.66'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.43%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.18%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3822'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -2.91%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3419'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -1.19%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '47.02'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.55%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.147'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -4.11%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07375'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.46%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '23.26'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +6.75%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.001'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.07%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.389'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.93%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.387'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +3.47%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.777.64'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.25%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001078'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.48%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.06650'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.55%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '82.43'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -2.74%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.001'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.22%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.39'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.82%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.400'
$ws.Range('D22').Style = "Normal"

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '28.157.08'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +2.87%  '

# Row 24
$ws.Range('E24').Value = '  -3.00%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.387'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.99%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.446'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.68%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '20.75'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -2.22%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.404'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -3.86%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '153.71'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -3.00%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.978.12'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.43%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '134.67'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.93%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.038'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.20%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.056'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +0.73%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.08937'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +1.03%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '12.70'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.31%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.02407'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.64%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.6836'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.24%  '

# Row 38
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.06374'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -2.29%  '

# Row 39
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.317'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.97%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.2159'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.42%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.240'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.85%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.494'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -8.10%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.181'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -2.06%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '14.29'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.47%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.000'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +0.13%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.6266'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.80%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.861'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.19%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '133.02'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.60%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.070'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.99%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.07516'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +4.83%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.205'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +2.83%  '
